$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E. This shifts the old "Rejected"
# column (and all its data) from E -> F, and copies formatting from the
# neighbouring column (D) onto the freshly inserted column E.
$ws.Columns("E").Insert()

# Give the new column its header and width (Status column).
$ws.Range("E1").Value = "Status"
$ws.Columns("E").ColumnWidth = 9.3

# Move the active cell/selection to the new Status column.
$ws.Range("E2").Select()
